# Updated graph and excel files
#
# Adds a chart title ("Distance calibration"), axis titles for the X axis
# ("Real distance (mm)") and Y axis ("Measured distance (mm)"),
# resizes/repositions the chart on the sheet, and moves the active cell
# selection to B8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlValue / xlCategory axis-type constants used by Chart.Axes()
$xlCategory = 1
$xlValue = 2
$xlPrimary = 1

$co = $ws.ChartObjects(1)
$chart = $co.Chart

# --- Chart title -----------------------------------------------------
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Distance calibration"

# --- Axis titles -------------------------------------------------------
$xAxis = $chart.Axes($xlCategory, $xlPrimary)
$xAxis.HasTitle = $true
$xAxis.AxisTitle.Text = "Real distance (mm)"

$yAxis = $chart.Axes($xlValue, $xlPrimary)
$yAxis.HasTitle = $true
$yAxis.AxisTitle.Text = "Measured distance (mm)"

# --- Move / resize the chart on the worksheet ---------------------------
# Old anchor: from col 14 (+238125 EMU), row 4 (+123825 EMU)
#               to col 21 (+542925 EMU), row 19 (+9525 EMU)
# New anchor: from col 3 (+276224 EMU), row 0 (+104775 EMU)
#               to col 24 (+171450 EMU), row 35 (+123825 EMU)
$co.Left = 197.06242125984252
$co.Top = 8.25
$co.Width = 1218.9375787401575
$co.Height = 526.5

# --- Worksheet selection -------------------------------------------------
$ws.Range("B8").Select()
